# Refresh "想去人数" (interest count) figures in column F across all sheets,
# matching the regenerated data snapshot (gh-pages output at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(7, 6).Value = 12747
$ws.Cells.Item(10, 6).Value = 2838
$ws.Cells.Item(11, 6).Value = 78
$ws.Cells.Item(12, 6).Value = 6089
$ws.Cells.Item(15, 6).Value = 3225
$ws.Cells.Item(16, 6).Value = 21
$ws.Cells.Item(17, 6).Value = 152
$ws.Cells.Item(20, 6).Value = 28
$ws.Cells.Item(23, 6).Value = 3503
$ws.Cells.Item(26, 6).Value = 2616
$ws.Cells.Item(28, 6).Value = 1834
$ws.Cells.Item(30, 6).Value = 185
$ws.Cells.Item(31, 6).Value = 6373
$ws.Cells.Item(32, 6).Value = 11
$ws.Cells.Item(33, 6).Value = 156
$ws.Cells.Item(34, 6).Value = 1921
$ws.Cells.Item(35, 6).Value = 1282
$ws.Cells.Item(40, 6).Value = 12
$ws.Cells.Item(41, 6).Value = 210
$ws.Cells.Item(44, 6).Value = 112
$ws.Cells.Item(45, 6).Value = 1153
$ws.Cells.Item(46, 6).Value = 1695
$ws.Cells.Item(49, 6).Value = 1154
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(10, 6).Value = 118
$ws.Cells.Item(14, 6).Value = 923
$ws.Cells.Item(16, 6).Value = 89
$ws.Cells.Item(19, 6).Value = 6
$ws.Cells.Item(20, 6).Value = 9
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 394
$ws.Cells.Item(3, 6).Value = 558
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 6).Value = 394
$ws.Cells.Item(9, 6).Value = 12747
$ws.Cells.Item(12, 6).Value = 2838
$ws.Cells.Item(13, 6).Value = 78
$ws.Cells.Item(14, 6).Value = 6089
$ws.Cells.Item(16, 6).Value = 3225
$ws.Cells.Item(17, 6).Value = 21
$ws.Cells.Item(18, 6).Value = 152
$ws.Cells.Item(20, 6).Value = 28
$ws.Cells.Item(22, 6).Value = 118
$ws.Cells.Item(24, 6).Value = 3503
$ws.Cells.Item(26, 6).Value = 2616
$ws.Cells.Item(28, 6).Value = 1834
$ws.Cells.Item(30, 6).Value = 185
$ws.Cells.Item(31, 6).Value = 6373
$ws.Cells.Item(32, 6).Value = 89
$ws.Cells.Item(33, 6).Value = 11
$ws.Cells.Item(34, 6).Value = 156
$ws.Cells.Item(36, 6).Value = 1282
$ws.Cells.Item(40, 6).Value = 210
$ws.Cells.Item(43, 6).Value = 113
$ws.Cells.Item(44, 6).Value = 1153
$ws.Cells.Item(45, 6).Value = 9
$ws.Cells.Item(46, 6).Value = 1695
$ws.Cells.Item(49, 6).Value = 1154
